$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump up the weight values (column Q) and multiplier (column P) for a few items,
# per "a lil more weight".

# Steyr AUG A1 Swarovski Optik 1.5x Scoped (row 5): weight 0.59 -> 0.6
$ws.Range("Q5").Value = 0.6

# Steyr AUG A3SF 60mm 3x scope (row 18): count 1 -> 2, weight 0.23 -> 0.26
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 0.26

# Steyr AUG A3SF 60mm 1.5x scope (row 19): count 1 -> 2, weight 0.23 -> 0.26
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 0.26

# Leave the active selection on Q19, matching the last-edited cell.
$ws.Range("Q19").Select()
